$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated activity values (runs, balls, fours, sixes). Values are kept as
# text (leading apostrophe) to match the workbook's existing convention of
# storing these numeric-looking stats as text cells.
$ws.Range("C2").Value = "'8"
$ws.Range("D2").Value = "'5"
$ws.Range("E2").Value = "'0"

$ws.Range("C3").Value = "'31"
$ws.Range("D3").Value = "'24"
$ws.Range("E3").Value = "'4"

$ws.Range("C5").Value = "'10"
$ws.Range("D5").Value = "'13"

$ws.Range("C6").Value = "'1"
$ws.Range("D6").Value = "'1"
$ws.Range("E6").Value = "'0"
$ws.Range("F6").Value = "'0"

$ws.Range("C7").Value = "'12"
$ws.Range("D7").Value = "'6"
$ws.Range("E7").Value = "'1"
$ws.Range("F7").Value = "'1"
